$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AB2").Value = 0.4399577731552838
$ws.Range("AB3").Value = 0.4208339345754986
$ws.Range("AB4").Value = 0.2751194273529608
$ws.Range("AB5").Value = 0.3426384194471387
$ws.Range("AB6").Value = 0.4052211072857823
$ws.Range("AB7").Value = 0.4397452894942148
